$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 580.59375
$ws.Range("I15").Value = 580.59375
$ws.Range("K15").Value = 1741.78125
$ws.Range("M15").Value = -1572.78125
$ws.Range("H43").Value = 55560532
$ws.Range("I43").Value = 100001064
$ws.Range("K43").Value = 100001064
$ws.Range("M43").Value = -100000995
$ws.Range("H55").Value = 221.76923
$ws.Range("I55").Value = 215.375
$ws.Range("K55").Value = 215.375
$ws.Range("M55").Value = -1.375
$ws.Range("H132").Value = 37041828
$ws.Range("I132").Value = 43483364
$ws.Range("K132").Value = 130450092
$ws.Range("M132").Value = -130447562
$ws.Range("H137").Value = 49785.785
$ws.Range("I137").Value = 82056.55
$ws.Range("J137").Value = 2455.3333
$ws.Range("K137").Value = 246169.65
$ws.Range("L137").Value = 7365.999899999999
$ws.Range("M137").Value = -243619.65
$ws.Range("N137").Value = -12465.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 21632.307
$ws.Range("I74").Value = 3285.973
$ws.Range("K74").Value = 3285.973
$ws.Range("M74").Value = -2411.973
$ws.Range("H77").Value = 21632.307
$ws.Range("I77").Value = 3285.973
$ws.Range("K77").Value = 16429.865
$ws.Range("M77").Value = -12061.865
$ws.Range("H88").Value = 1682.625
$ws.Range("J88").Value = 1135.3334
$ws.Range("L88").Value = 1135.3334
$ws.Range("N88").Value = -1947.3334
$ws.Range("H91").Value = 1682.625
$ws.Range("J91").Value = 1135.3334
$ws.Range("L91").Value = 1135.3334
$ws.Range("N91").Value = -3943.3334
$ws.Range("H97").Value = 3236403.5
$ws.Range("I97").Value = 4044259.8
$ws.Range("J97").Value = 4978.5
$ws.Range("K97").Value = 4044259.8
$ws.Range("L97").Value = 4978.5
$ws.Range("M97").Value = -4043763.8
$ws.Range("N97").Value = -5970.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1846
$ws.Range("I7").Value = 1215.2
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 1215.2
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -1102.2
$ws.Range("N7").Value = -5226
$ws.Range("H26").Value = 16231
$ws.Range("I26").Value = 16231
$ws.Range("K26").Value = 16231
$ws.Range("M26").Value = -15939
$ws.Range("H96").Value = 22987.834
$ws.Range("I96").Value = 22987.834
$ws.Range("K96").Value = 22987.834
$ws.Range("M96").Value = -20241.834
$ws.Range("H139").Value = 163738.12
$ws.Range("J139").Value = 174151
$ws.Range("L139").Value = 174151
$ws.Range("N139").Value = -184431

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1224
$ws.Range("H31").Value = 19662.775
$ws.Range("I31").Value = 2767.56
$ws.Range("J31").Value = 32462.182
$ws.Range("K31").Value = 2767.56
$ws.Range("L31").Value = 32462.182
$ws.Range("M31").Value = -2472.56
$ws.Range("N31").Value = -33052.182
$ws.Range("H34").Value = 19662.775
$ws.Range("I34").Value = 2767.56
$ws.Range("J34").Value = 32462.182
$ws.Range("K34").Value = 2767.56
$ws.Range("L34").Value = 32462.182
$ws.Range("M34").Value = -2565.56
$ws.Range("N34").Value = -32866.182
$ws.Range("H58").Value = 7563.125
$ws.Range("I58").Value = 9589.071
$ws.Range("J58").Value = 4726.8
$ws.Range("K58").Value = 9589.071
$ws.Range("L58").Value = 4726.8
$ws.Range("M58").Value = -9386.071
$ws.Range("N58").Value = -5132.8
$ws.Range("H132").Value = 97070.89999999999
$ws.Range("I132").Value = 78660.766
$ws.Range("K132").Value = 235982.298
$ws.Range("M132").Value = -233452.298
$ws.Range("H134").Value = 3764
$ws.Range("I134").Value = 2127.4546
$ws.Range("J134").Value = 5400.5454
$ws.Range("K134").Value = 6382.3638
$ws.Range("L134").Value = 16201.6362
$ws.Range("M134").Value = -3847.3638
$ws.Range("N134").Value = -21271.6362
$ws.Range("H136").Value = 7563.125
$ws.Range("I136").Value = 9589.071
$ws.Range("J136").Value = 4726.8
$ws.Range("K136").Value = 28767.213
$ws.Range("L136").Value = 14180.4
$ws.Range("M136").Value = -26217.213
$ws.Range("N136").Value = -19280.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1768.0588
$ws.Range("I3").Value = 1196.6923
$ws.Range("K3").Value = 3590.0769
$ws.Range("M3").Value = -3478.0769
$ws.Range("H103").Value = 281.2857
$ws.Range("I103").Value = 281.2857
$ws.Range("K103").Value = 843.8571000000001
$ws.Range("M103").Value = 35.14289999999994
$ws.Range("H108").Value = 955
$ws.Range("I108").Value = 955
$ws.Range("K108").Value = 2865
$ws.Range("M108").Value = 15
$ws.Range("H133").Value = 3226.3333
$ws.Range("I133").Value = 951
$ws.Range("J133").Value = 7777
$ws.Range("K133").Value = 2853
$ws.Range("L133").Value = 23331
$ws.Range("M133").Value = 2207
$ws.Range("N133").Value = -33451

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3427.4348
$ws.Range("I132").Value = 3356.182
$ws.Range("J132").Value = 4995
$ws.Range("K132").Value = 10068.546
$ws.Range("L132").Value = 14985
$ws.Range("M132").Value = -7538.545999999998
$ws.Range("N132").Value = -20045

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10154.6
$ws.Range("I7").Value = 8181.75
$ws.Range("J7").Value = 10872
$ws.Range("K7").Value = 8181.75
$ws.Range("L7").Value = 10872
$ws.Range("M7").Value = -8069.75
$ws.Range("N7").Value = -11096
$ws.Range("H96").Value = 69741.5
$ws.Range("J96").Value = 69741.5
$ws.Range("L96").Value = 69741.5
$ws.Range("N96").Value = -75233.5
$ws.Range("H126").Value = 10154.6
$ws.Range("I126").Value = 8181.75
$ws.Range("J126").Value = 10872
$ws.Range("K126").Value = 24545.25
$ws.Range("L126").Value = 32616
$ws.Range("M126").Value = -22075.25
$ws.Range("N126").Value = -37556
$ws.Range("H136").Value = 36374.54
$ws.Range("I136").Value = 49892.74
$ws.Range("K136").Value = 149678.22
$ws.Range("M136").Value = -147128.22

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 10010000
$ws.Range("I11").Value = 20000000
$ws.Range("K11").Value = 20000000
$ws.Range("M11").Value = -19999858
$ws.Range("H62").Value = 7656.8887
$ws.Range("J62").Value = 7656.8887
$ws.Range("L62").Value = 7656.8887
$ws.Range("N62").Value = -8904.8887
$ws.Range("H65").Value = 7656.8887
$ws.Range("J65").Value = 7656.8887
$ws.Range("L65").Value = 38284.4435
$ws.Range("N65").Value = -44524.4435
$ws.Range("H105").Value = 33638.332
$ws.Range("J105").Value = 33638.332
$ws.Range("L105").Value = 33638.332
$ws.Range("N105").Value = -40626.332
$ws.Range("H107").Value = 200002430
$ws.Range("I107").Value = 250002260
$ws.Range("K107").Value = 750006780
$ws.Range("M107").Value = -750004860
$ws.Range("H113").Value = 1537.4814
$ws.Range("I113").Value = 713.2941
$ws.Range("J113").Value = 2938.6
$ws.Range("K113").Value = 2139.8823
$ws.Range("L113").Value = 8815.799999999999
$ws.Range("M113").Value = 30.11770000000024
$ws.Range("N113").Value = -13155.8
$ws.Range("H122").Value = 2912.3635
$ws.Range("I122").Value = 2813.1904
$ws.Range("J122").Value = 4995
$ws.Range("K122").Value = 8439.5712
$ws.Range("L122").Value = 14985
$ws.Range("M122").Value = -5989.5712
$ws.Range("N122").Value = -19885
